$d = $word.ActiveDocument
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq [char]13 -and $p.Next().Range.Text.TrimEnd([char]13,[char]7) -eq "Tools") {
        $target = $p
        break
    }
}
if ($target -eq $null) { throw "target paragraph not found" }
$xmlFragment = '<w:p><w:pPr><w:pStyle w:val="Default"/><w:spacing w:line="280" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Default"/><w:spacing w:line="280" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>Search</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Default"/><w:spacing w:line="280" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t>To setup search, we originally tried going with Google</w:t></w:r><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>’</w:t></w:r><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">s Custom Search Engine. Since we were generating pages dynamically, we had to provide the CSE with a sitemap and give it a couple of days to process the information on our pages so that it would be able to display results. However, we ran into a few issues. The results were sparse and the search results were extremely slow, ranging anywhere from five to fifteen seconds. We decided that this search would be a negative experience on the user and drive them away from the platform, rather than attract new viewers and provide useful functionality. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Default"/><w:spacing w:line="280" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t>On the frontend side, we created a component for the search page and captured the fields that the user was searching for from the search bar. We then dynamically routed their search query to the url /search/query and displayed those results. We displayed the results from our search similar to a Google search, where the title of the entry is the link to its page, and all of its relevant information is right below it for the user to easily access and read. This functionality of displaying contextualization is a big benefit to the target users of our platform, who do not need to go through and check out every entry one at a time, page by page, to absorb information. This allows them to easily process desired information or skim over entries that aren</w:t></w:r><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>’</w:t></w:r><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">t relevant to their search. We also created a default no results found page to indicate to a user that their query is bad. It could be the case that the user may think that our database simply does not have the information they are searching for, and this case prompted us to distinguish those two scenarios so that the user can be more well informed. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Default"/><w:spacing w:line="280" w:lineRule="atLeast"/></w:pPr><w:r><w:rPr><w:rStyle w:val="None"/><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">The frontend simply makes a call to our API which searches for the captured query from the user. The backend took this query, and searched through the database using the query as a regex. This search was through all attributes of all the entries of all of our models, but did not include foreign key attributes (which were simply ids or meaningless numbers to the user). The backend then returned a nicely packed JSON, which essentially was a list of lists of models. The first list indexed into four separate other lists, one for each of our models, and those four separate other lists contained all the entries that contained values matching the users search query. The frontend then parsed this data accordingly. Then, four new components were created as a means of displaying the results, one for each type of model. The components were then created dynamically, as one entry was parsed it was sent to the components which handled that information and visualization. We used a different component for each type of model because our models all had their own unique attributes. </w:t></w:r></w:p>'
$target.Range.InsertXML($xmlFragment)
